# DataOrangeHrmLive.xlsx - add "Recruitment" sheet/data, per commit:
# "Creacion de recruitment.feature y libro recruitment"

$wb = $excel.ActiveWorkbook

# --- Login sheet: restore cursor to A1 before losing focus -----------------
$login = $wb.Worksheets.Item("Login")
$login.Range("A1").Select() | Out-Null

# --- Create the new "Recruitment" sheet right after "Login" ----------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $login)
$ws.Name = "Recruitment"

# Reuse Login's bold/centered header style (cellXf index 1) for row 1,
# A1:M1, by pasting formats from Login!A1:C1 (tiles across the range).
$login.Range("A1:C1").Copy()
$ws.Range("A1:M1").PasteSpecial(-4122)

# --- Header row --------------------------------------------------------
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "FIRST_NAME"
$ws.Range("C1").Value = "MIDDLE_NAME"
$ws.Range("D1").Value = "LAST_NAME"
$ws.Range("E1").Value = "VACANCY"
$ws.Range("F1").Value = "EMAIL"
$ws.Range("G1").Value = "CONTACT_NUMBER"
$ws.Range("H1").Value = "KEYWORDS"
$ws.Range("I1").Value = "DATA_OF_APPLICATION"
$ws.Range("J1").Value = "NOTE"
$ws.Range("K1").Value = "CONSENT_TO_KEEP_DATA"
$ws.Range("L1").Value = "SHORTLIST_NOTE"

# --- Data row ------------------------------------------------------------
$ws.Range("B2").Value = "Wilder "
$ws.Range("C2").Value = "de Jesus "
$ws.Range("D2").Value = "Bernal Lopez"
$ws.Range("E2").Value = "QA LEAD"
$ws.Range("F2").Value = "w@w.com"
$ws.Range("G2").Value = 3126148527
$ws.Range("H2").Value = "El mejor"
$ws.Range("I2").Value = "2024-22-07"
$ws.Range("J2").Value = "Nota1"
$ws.Range("K2").Value = $true
$ws.Range("K2").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("L2").Value = "Nota2"
$ws.Rows.Item(2).RowHeight = 15.65

# Hyperlink on the e-mail cell (Hyperlinks.Add forces the built-in
# "Hyperlink" style onto the cell as a side effect, so strip that back off
# the cell afterwards to keep it on the plain/default style).
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:w@w.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "w@w.com") | Out-Null
$ws.Range("F2").ClearFormats()
$ws.Range("F2").Value = "w@w.com"

# --- Column widths (characters, approximating the authored widths) ------
$ws.Columns.Item(1).ColumnWidth = 2.026667
$ws.Columns.Item(2).ColumnWidth = 10.386667
$ws.Columns.Item(3).ColumnWidth = 12.006667
$ws.Columns.Item(4).ColumnWidth = 10.696667
$ws.Columns.Item(5).ColumnWidth = 7.766667
$ws.Columns.Item(6).ColumnWidth = 8.576667
$ws.Columns.Item(7).ColumnWidth = 15.336667
$ws.Columns.Item(8).ColumnWidth = 9.476667
$ws.Columns.Item(9).ColumnWidth = 18.566667
$ws.Columns.Item(11).ColumnWidth = 20.576667
$ws.Columns.Item(12).ColumnWidth = 14.326667

# --- Page setup / margins / header-footer --------------------------------
$ws.PageSetup.LeftMargin = 56.7
$ws.PageSetup.RightMargin = 56.7
$ws.PageSetup.TopMargin = 75.8
$ws.PageSetup.BottomMargin = 75.8
$ws.PageSetup.HeaderMargin = 56.7
$ws.PageSetup.FooterMargin = 56.7
$ws.PageSetup.PaperSize = 1
$ws.PageSetup.Zoom = 100
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.Order = 1
$ws.PageSetup.Orientation = 1
$ws.PageSetup.BlackAndWhite = $false
$ws.PageSetup.Draft = $false
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&Kffffff&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&KffffffPage &P'

# --- Selection left where the author last worked -------------------------
$ws.Range("K10").Select() | Out-Null

Write-Output "Recruitment sheet added"
